# Refined metadata to be additional tab
#
# 1) Update the existing "data" sheet's time_taken (F) column timestamps.
# 2) Add a new "metadata" sheet (placed after "data") describing the
#    PanelApp query that produced this workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = "2021-10-05 14:22:14.686148"
$ws1.Range("F3").Value = "2021-10-05 14:22:14.686157"
$ws1.Range("F4").Value = "2021-10-05 14:22:14.686161"
$ws1.Range("F5").Value = "2021-10-05 14:22:14.686164"
$ws1.Range("F6").Value = "2021-10-05 14:22:14.686167"
$ws1.Range("F7").Value = "2021-10-05 14:22:14.686170"
$ws1.Range("F8").Value = "2021-10-05 14:22:14.686173"
$ws1.Range("F9").Value = "2021-10-05 14:22:14.686175"
$ws1.Range("F10").Value = "2021-10-05 14:22:14.686178"
$ws1.Range("F11").Value = "2021-10-05 14:22:14.686181"
$ws1.Range("F12").Value = "2021-10-05 14:22:14.686184"
$ws1.Range("F13").Value = "2021-10-05 14:22:14.686186"
$ws1.Range("F14").Value = "2021-10-05 14:22:14.686189"
$ws1.Range("F15").Value = "2021-10-05 14:22:14.686192"
$ws1.Range("F16").Value = "2021-10-05 14:22:14.686194"
$ws1.Range("F17").Value = "2021-10-05 14:22:14.686197"
$ws1.Range("F18").Value = "2021-10-05 14:22:14.686200"
$ws1.Range("F19").Value = "2021-10-05 14:22:14.686202"
$ws1.Range("F20").Value = "2021-10-05 14:22:14.686205"
$ws1.Range("F21").Value = "2021-10-05 14:22:14.686208"
$ws1.Range("F22").Value = "2021-10-05 14:22:14.686211"
$ws1.Range("F23").Value = "2021-10-05 14:22:14.686214"
$ws1.Range("F24").Value = "2021-10-05 14:22:14.686216"
$ws1.Range("F25").Value = "2021-10-05 14:22:14.686219"
$ws1.Range("F26").Value = "2021-10-05 14:22:14.686222"
$ws1.Range("F27").Value = "2021-10-05 14:22:14.686225"
$ws1.Range("F28").Value = "2021-10-05 14:22:14.686228"
$ws1.Range("F29").Value = "2021-10-05 14:22:14.686231"
$ws1.Range("F30").Value = "2021-10-05 14:22:14.686233"
$ws1.Range("F31").Value = "2021-10-05 14:22:14.686236"
$ws1.Range("F32").Value = "2021-10-05 14:22:14.686239"
$ws1.Range("F33").Value = "2021-10-05 14:22:14.686241"
$ws1.Range("F34").Value = "2021-10-05 14:22:14.686245"
$ws1.Range("F35").Value = "2021-10-05 14:22:14.686247"
$ws1.Range("F36").Value = "2021-10-05 14:22:14.686250"

# --- Add the "metadata" sheet, positioned immediately after "data" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "metadata"

# Reuse the bold/bordered header style from the "data" sheet's header row
# (and the index-column style from A2) instead of defining new styles.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws2.Range("F1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Pneumothorax - familial"
$ws2.Range("C2").Value = 105
$ws2.Range("D2").Value = "'2.37"
$ws2.Range("D2").Style = "Normal"
$ws2.Range("E2").Value = "2021-03-17T15:04:38.037565Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:14.682592"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/105/?format=json"

$ws1.Select()
$ws1.Range("A1").Select()
